# Append new ticker rows to the sheet (data update 2024-03-05)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTickers = @("IMX-USD", "TAO-USD", "GRT-USD", "MNT-USD", "PEPE-USD")

# Find the last used row in column A and append right after it.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 0 }

$startRow = $lastRow + 1
for ($i = 0; $i -lt $newTickers.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newTickers[$i]
}
